{"js": "// Replace each old expression/date with its new value, in document order.\n// Each (oldText -> newText) pair corresponds to one run of text in the\n// title paragraph (the date) or one table cell (the math expression).\nconst replacements = [\n  [\"2025-02-21 Friday\", \"2025-02-22 Saturday\"],\n  [\"5+31=36\", \"55-7=48\"],\n  [\"72+7=79\", \"9+7=16\"],\n  [\"71+13=84\", \"43-7=36\"],\n  [\"68-24=44\", \"38+38=76\"],\n  [\"39+12=51\", \"75-35=40\"],\n  [\"81-18=63\", \"33+2=35\"],\n  [\"51-4=47\", \"57+8=65\"],\n  [\"0+34=34\", \"64+24=88\"],\n  [\"8+48=56\", \"56-33=23\"],\n  [\"23+72=95\", \"41-40=1\"],\n  [\"4+92=96\", \"46-4=42\"],\n  [\"5+83=88\", \"0+92=92\"],\n  [\"4+14=18\", \"80-19=61\"],\n  [\"26-14=12\", \"9+65=74\"],\n  [\"77-44=33\", \"95-13=82\"],\n  [\"13+0=13\", \"63-58=5\"],\n  [\"28+38=66\", \"27+8=35\"],\n  [\"16+48=64\", \"37+3=40\"],\n  [\"18-15=3\", \"90-76=14\"],\n  [\"26+62=88\", \"90-88=2\"],\n  [\"66-29=37\", \"52-14=38\"],\n  [\"1+23=24\", \"8+24=32\"],\n  [\"21+63=84\", \"4+58=62\"],\n  [\"28-21=7\", \"20-18=2\"],\n  [\"88-11=77\", \"33+15=48\"],\n  [\"66+26=92\", \"68+8=76\"],\n  [\"37-5=32\", \"36-25=11\"],\n  [\"54+11=65\", \"77-60=17\"],\n  [\"16+61=77\", \"38+3=41\"],\n  [\"91-12=79\", \"98-86=12\"],\n  [\"60-40=20\", \"39+47=86\"],\n  [\"59+12=71\", \"75-26=49\"],\n  [\"53-43=10\", \"48+22=70\"],\n  [\"81-9=72\", \"26+55=81\"],\n  [\"91-73=18\", \"28+31=59\"],\n  [\"29+29=58\", \"2+25=27\"],\n  [\"55+30=85\", \"97-42=55\"],\n  [\"35+10=45\", \"24+17=41\"],\n  [\"80-77=3\", \"93-19=74\"],\n  [\"91-78=13\", \"91-68=23\"],\n  [\"29+26=55\", \"36+38=74\"],\n  [\"34+11=45\", \"17+75=92\"],\n  [\"2+96=98\", \"76-49=27\"],\n  [\"1+14=15\", \"91-64=27\"],\n  [\"87+6=93\", \"26+14=40\"],\n  [\"8-5=3\", \"51+41=92\"],\n  [\"85+13=98\", \"62+18=80\"],\n  [\"93-7=86\", \"73-22=51\"],\n  [\"52-3=49\", \"24+43=67\"],\n  [\"21+70=91\", \"85-26=59\"],\n  [\"7+50=57\", \"82-26=56\"],\n  [\"16+23=39\", \"16-10=6\"],\n  [\"47+20=67\", \"24-5=19\"],\n  [\"31-15=16\", \"2+79=81\"],\n  [\"12+34=46\", \"89-0=89\"],\n  [\"34-6=28\", \"69-14=55\"],\n  [\"89-26=63\", \"57+4=61\"],\n  [\"96-3=93\", \"77-42=35\"],\n  [\"44-34=10\", \"63-3=60\"],\n  [\"30+7=37\", \"29-10=19\"],\n  [\"80-76=4\", \"90-88=2\"],\n  [\"61-58=3\", \"79-25=54\"],\n  [\"52-31=21\", \"77-10=67\"],\n  [\"16+71=87\", \"92-40=52\"],\n  [\"10+68=78\", \"66+33=99\"],\n  [\"32+21=53\", \"55+18=73\"],\n  [\"17+17=34\", \"3+96=99\"],\n  [\"76-42=34\", \"86-55=31\"],\n  [\"64-12=52\", \"66-62=4\"],\n  [\"17+10=27\", \"44+52=96\"],\n  [\"94-19=75\", \"40+46=86\"],\n  [\"48-23=25\", \"61+35=96\"],\n  [\"23+9=32\", \"33+46=79\"],\n  [\"61-34=27\", \"91-29=62\"],\n  [\"39-3=36\", \"3+18=21\"],\n  [\"74-48=26\", \"16+75=91\"],\n  [\"87+11=98\", \"5+90=95\"],\n  [\"94-12=82\", \"50-33=17\"],\n  [\"29+59=88\", \"24-4=20\"],\n  [\"7+43=50\", \"90-87=3\"],\n  [\"36+34=70\", \"13-1=12\"],\n  [\"35+49=84\", \"32-10=22\"],\n  [\"23+56=79\", \"17+3=20\"],\n  [\"88-3=85\", \"47+50=97\"],\n  [\"76+2=78\", \"96-53=43\"],\n  [\"48+51=99\", \"36-31=5\"],\n  [\"63-1=62\", \"20-19=1\"],\n  [\"61-23=38\", \"40+5=45\"],\n  [\"73-31=42\", \"51-21=30\"],\n  [\"29+45=74\", \"51-3=48\"],\n  [\"91-42=49\", \"31+46=77\"],\n  [\"75-22=53\", \"3+87=90\"],\n  [\"13+14=27\", \"96-38=58\"],\n  [\"56+20=76\", \"12-9=3\"],\n  [\"27+3=30\", \"35-7=28\"],\n  [\"72-5=67\", \"45+44=89\"],\n  [\"66-40=26\", \"33-28=5\"],\n  [\"52-6=46\", \"61-25=36\"],\n  [\"89-5=84\", \"81-23=58\"],\n  [\"49-16=33\", \"1+53=54\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  // Search the whole body for the exact (case-sensitive) old text.\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  // Replace just that range's text in place (keeps existing run formatting).\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace each old expression/date with its new value, in document order.\n# Each pair corresponds to one run of text in the title paragraph (the\n# date) or one table cell (the math expression).\n$pairs = @(\n  @('2025-02-21 Friday', '2025-02-22 Saturday'),\n  @('5+31=36', '55-7=48'),\n  @('72+7=79', '9+7=16'),\n  @('71+13=84', '43-7=36'),\n  @('68-24=44', '38+38=76'),\n  @('39+12=51', '75-35=40'),\n  @('81-18=63', '33+2=35'),\n  @('51-4=47', '57+8=65'),\n  @('0+34=34', '64+24=88'),\n  @('8+48=56', '56-33=23'),\n  @('23+72=95', '41-40=1'),\n  @('4+92=96', '46-4=42'),\n  @('5+83=88', '0+92=92'),\n  @('4+14=18', '80-19=61'),\n  @('26-14=12', '9+65=74'),\n  @('77-44=33', '95-13=82'),\n  @('13+0=13', '63-58=5'),\n  @('28+38=66', '27+8=35'),\n  @('16+48=64', '37+3=40'),\n  @('18-15=3', '90-76=14'),\n  @('26+62=88', '90-88=2'),\n  @('66-29=37', '52-14=38'),\n  @('1+23=24', '8+24=32'),\n  @('21+63=84', '4+58=62'),\n  @('28-21=7', '20-18=2'),\n  @('88-11=77', '33+15=48'),\n  @('66+26=92', '68+8=76'),\n  @('37-5=32', '36-25=11'),\n  @('54+11=65', '77-60=17'),\n  @('16+61=77', '38+3=41'),\n  @('91-12=79', '98-86=12'),\n  @('60-40=20', '39+47=86'),\n  @('59+12=71', '75-26=49'),\n  @('53-43=10', '48+22=70'),\n  @('81-9=72', '26+55=81'),\n  @('91-73=18', '28+31=59'),\n  @('29+29=58', '2+25=27'),\n  @('55+30=85', '97-42=55'),\n  @('35+10=45', '24+17=41'),\n  @('80-77=3', '93-19=74'),\n  @('91-78=13', '91-68=23'),\n  @('29+26=55', '36+38=74'),\n  @('34+11=45', '17+75=92'),\n  @('2+96=98', '76-49=27'),\n  @('1+14=15', '91-64=27'),\n  @('87+6=93', '26+14=40'),\n  @('8-5=3', '51+41=92'),\n  @('85+13=98', '62+18=80'),\n  @('93-7=86', '73-22=51'),\n  @('52-3=49', '24+43=67'),\n  @('21+70=91', '85-26=59'),\n  @('7+50=57', '82-26=56'),\n  @('16+23=39', '16-10=6'),\n  @('47+20=67', '24-5=19'),\n  @('31-15=16', '2+79=81'),\n  @('12+34=46', '89-0=89'),\n  @('34-6=28', '69-14=55'),\n  @('89-26=63', '57+4=61'),\n  @('96-3=93', '77-42=35'),\n  @('44-34=10', '63-3=60'),\n  @('30+7=37', '29-10=19'),\n  @('80-76=4', '90-88=2'),\n  @('61-58=3', '79-25=54'),\n  @('52-31=21', '77-10=67'),\n  @('16+71=87', '92-40=52'),\n  @('10+68=78', '66+33=99'),\n  @('32+21=53', '55+18=73'),\n  @('17+17=34', '3+96=99'),\n  @('76-42=34', '86-55=31'),\n  @('64-12=52', '66-62=4'),\n  @('17+10=27', '44+52=96'),\n  @('94-19=75', '40+46=86'),\n  @('48-23=25', '61+35=96'),\n  @('23+9=32', '33+46=79'),\n  @('61-34=27', '91-29=62'),\n  @('39-3=36', '3+18=21'),\n  @('74-48=26', '16+75=91'),\n  @('87+11=98', '5+90=95'),\n  @('94-12=82', '50-33=17'),\n  @('29+59=88', '24-4=20'),\n  @('7+43=50', '90-87=3'),\n  @('36+34=70', '13-1=12'),\n  @('35+49=84', '32-10=22'),\n  @('23+56=79', '17+3=20'),\n  @('88-3=85', '47+50=97'),\n  @('76+2=78', '96-53=43'),\n  @('48+51=99', '36-31=5'),\n  @('63-1=62', '20-19=1'),\n  @('61-23=38', '40+5=45'),\n  @('73-31=42', '51-21=30'),\n  @('29+45=74', '51-3=48'),\n  @('91-42=49', '31+46=77'),\n  @('75-22=53', '3+87=90'),\n  @('13+14=27', '96-38=58'),\n  @('56+20=76', '12-9=3'),\n  @('27+3=30', '35-7=28'),\n  @('72-5=67', '45+44=89'),\n  @('66-40=26', '33-28=5'),\n  @('52-6=46', '61-25=36'),\n  @('89-5=84', '81-23=58'),\n  @('49-16=33', '1+53=54')\n)\n\n$d = $word.ActiveDocument\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n  #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n  #   ReplaceWith, Replace)\n  # MatchCase=$true (exact text), MatchWildcards=$false (literal \"+\"/\"-\"),\n  # Replace=2 -> wdReplaceAll (each old value is unique, so this is a\n  # single in-place replacement keeping the run's existing formatting).\n  $found = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n  if (-not $found) {\n    Write-Output \"NOT FOUND: $old\"\n  }\n}\n"}
